$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 17 (Harry Potter and the Deathly Hallows)
# Use the raw date serial numbers (Excel 1900 date system) so no new
# number-format styles get created; then copy the date style from C17
# (style index 1, already mm/dd/yyyy) onto the new date cells.
$ws.Range("D17").Value = 44252
$ws.Range("C17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "fiction;wizards;adventure;harry potter"
$ws.Range("F17").Value = "Audio"
$ws.Range("G17").Value = "21 Hours 47 Mins"
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = $true

# Add row 18 (Think Again)
$ws.Range("A18").Value = "Think Again"
$ws.Range("B18").Value = "Adam Grant"
$ws.Range("C18").Value = 44250
$ws.Range("D18").Value = 44253
$ws.Range("C17").Copy()
$ws.Range("C18:D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "self improvement;business;rethinking"
$ws.Range("F18").Value = "Audio"
$ws.Range("G18").Value = "6 Hours 40 Mins"
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = $true

$excel.CutCopyMode = $false
$ws.Range("A19").Select()
